$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append one new quarterly data row (01-07-2021) as row 76, mirroring the
# existing "Serie" layout: column A is a text period label (shared string),
# columns B:BT are the numeric export figures for that quarter.
$rowNum = 76

# Write the period label as literal text (not an Excel date). A helper cell
# holds a formula whose result is the plain string "01-07-2021"; copying its
# *value* into the target cell preserves the text type instead of Excel's
# smart date re-interpretation of a typed "dd-mm-yyyy"-shaped string.
$helper = $ws.Cells.Item($ws.Rows.Count, 1)
$helper.Formula = "=""01-07-2021"""
$helper.Copy()
$ws.Cells.Item($rowNum, 1).PasteSpecial(-4163)
$helper.ClearContents()
$excel.CutCopyMode = $false

$values = @(
    23194, 14700, 13124, 5355, 7172, 839, 105, 203, 127, 252,
    34, 1129, 967, 3, 223, 24, 1, 103, 0, 0,
    24, 119, 10, 47, 13, 30, 7365, 2672, 121, 44,
    1102, 69, 27, 12, 173, 103, 116, 74, 31, 154,
    176, 529, 23, 405, 85, 722, 237, 69, 97, 108,
    135, 925, 86, 339, 271, 82, 1428, 35, 105, 98,
    135, 431, 83, 306, 90, 84, 519, 89, 261, 169,
    264
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 2
    $ws.Cells.Item($rowNum, $col).Value = $values[$i]
}
